$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 6).Value = 336
$ws.Cells.Item(2, 12).Value = 'stimuli/img_ac0ey.png'
$ws.Cells.Item(2, 13).Value = 86.62222222222222
$ws.Cells.Item(2, 14).Value = 70.02222222222223
$ws.Cells.Item(2, 15).Value = 78.32222222222222
$ws.Cells.Item(2, 16).Value = 45
$ws.Cells.Item(2, 17).Value = 9
$ws.Cells.Item(2, 18).Value = 9
$ws.Cells.Item(2, 19).Value = 9
$ws.Cells.Item(3, 6).Value = 337
$ws.Cells.Item(3, 8).Value = 'bedrooms'
$ws.Cells.Item(3, 9).Value = 'distractor'
$ws.Cells.Item(3, 11).Value = 'f'
$ws.Cells.Item(3, 12).Value = 'stimuli/img_e26ut.png'
$ws.Cells.Item(3, 13).Value = 81.07692307692308
$ws.Cells.Item(3, 14).Value = 61.28205128205128
$ws.Cells.Item(3, 15).Value = 71.17948717948718
$ws.Cells.Item(3, 16).Value = 39
$ws.Cells.Item(3, 17).Value = 8
$ws.Cells.Item(3, 18).Value = 8
$ws.Cells.Item(3, 19).Value = 8
$ws.Cells.Item(4, 6).Value = 338
$ws.Cells.Item(4, 12).Value = 'stimuli/img_91csq.png'
$ws.Cells.Item(4, 13).Value = 50.44736842105263
$ws.Cells.Item(4, 14).Value = 28.34210526315789
$ws.Cells.Item(4, 15).Value = 39.39473684210526
$ws.Cells.Item(4, 16).Value = 38
$ws.Cells.Item(5, 6).Value = 339
$ws.Cells.Item(5, 12).Value = 'stimuli/img_vgh2g.png'
$ws.Cells.Item(5, 13).Value = 93.81395348837209
$ws.Cells.Item(5, 14).Value = 78.27906976744185
$ws.Cells.Item(5, 15).Value = 86.04651162790697
$ws.Cells.Item(5, 16).Value = 43
$ws.Cells.Item(5, 17).Value = 10
$ws.Cells.Item(5, 18).Value = 10
$ws.Cells.Item(5, 19).Value = 10
$ws.Cells.Item(6, 6).Value = 340
$ws.Cells.Item(6, 8).Value = 'living_rooms'
$ws.Cells.Item(6, 9).Value = 'target'
$ws.Cells.Item(6, 11).Value = 'j'
$ws.Cells.Item(6, 12).Value = 'stimuli/img_g13d5.png'
$ws.Cells.Item(6, 13).Value = 73
$ws.Cells.Item(6, 14).Value = 51.51111111111111
$ws.Cells.Item(6, 15).Value = 62.25555555555556
$ws.Cells.Item(6, 16).Value = 45
$ws.Cells.Item(6, 17).Value = 6
$ws.Cells.Item(6, 18).Value = 6
$ws.Cells.Item(6, 19).Value = 6
$ws.Cells.Item(7, 6).Value = 341
$ws.Cells.Item(7, 12).Value = 'stimuli/img_73pyk.png'
$ws.Cells.Item(7, 13).Value = 69.27659574468085
$ws.Cells.Item(7, 14).Value = 47.27659574468085
$ws.Cells.Item(7, 15).Value = 58.27659574468085
$ws.Cells.Item(7, 16).Value = 47
$ws.Cells.Item(7, 17).Value = 5
$ws.Cells.Item(7, 18).Value = 5
$ws.Cells.Item(7, 19).Value = 5
$ws.Cells.Item(8, 6).Value = 342
$ws.Cells.Item(8, 8).Value = 'bedrooms'
$ws.Cells.Item(8, 9).Value = 'distractor'
$ws.Cells.Item(8, 11).Value = 'f'
$ws.Cells.Item(8, 12).Value = 'stimuli/img_mdpr4.png'
$ws.Cells.Item(8, 13).Value = 74.04255319148936
$ws.Cells.Item(8, 14).Value = 54.70212765957447
$ws.Cells.Item(8, 15).Value = 64.37234042553192
$ws.Cells.Item(8, 16).Value = 47
$ws.Cells.Item(8, 17).Value = 6
$ws.Cells.Item(8, 18).Value = 6
$ws.Cells.Item(8, 19).Value = 6
$ws.Cells.Item(9, 6).Value = 343
$ws.Cells.Item(9, 12).Value = 'stimuli/img_xr3up.png'
$ws.Cells.Item(9, 13).Value = 76.24444444444444
$ws.Cells.Item(9, 14).Value = 55.88888888888889
$ws.Cells.Item(9, 15).Value = 66.06666666666666
$ws.Cells.Item(9, 16).Value = 45
$ws.Cells.Item(9, 17).Value = 7
$ws.Cells.Item(9, 18).Value = 7
$ws.Cells.Item(9, 19).Value = 7
$ws.Cells.Item(10, 6).Value = 344
$ws.Cells.Item(10, 12).Value = 'stimuli/img_syam3.png'
$ws.Cells.Item(10, 13).Value = 41.32432432432432
$ws.Cells.Item(10, 14).Value = 26.2972972972973
$ws.Cells.Item(10, 15).Value = 33.81081081081081
$ws.Cells.Item(10, 16).Value = 37
$ws.Cells.Item(10, 17).Value = 2
$ws.Cells.Item(10, 18).Value = 2
$ws.Cells.Item(10, 19).Value = 2
$ws.Cells.Item(11, 6).Value = 345
$ws.Cells.Item(11, 8).Value = 'living_rooms'
$ws.Cells.Item(11, 9).Value = 'target'
$ws.Cells.Item(11, 11).Value = 'j'
$ws.Cells.Item(11, 12).Value = 'stimuli/img_tujn3.png'
$ws.Cells.Item(11, 13).Value = 81.40909090909091
$ws.Cells.Item(11, 14).Value = 62.52272727272727
$ws.Cells.Item(11, 15).Value = 71.96590909090909
$ws.Cells.Item(11, 16).Value = 44
$ws.Cells.Item(11, 17).Value = 8
$ws.Cells.Item(11, 18).Value = 8
$ws.Cells.Item(11, 19).Value = 8
$ws.Cells.Item(12, 6).Value = 346
$ws.Cells.Item(12, 8).Value = 'bedrooms'
$ws.Cells.Item(12, 9).Value = 'distractor'
$ws.Cells.Item(12, 11).Value = 'f'
$ws.Cells.Item(12, 12).Value = 'stimuli/img_gxm46.png'
$ws.Cells.Item(12, 13).Value = 74.78378378378379
$ws.Cells.Item(12, 14).Value = 54
$ws.Cells.Item(12, 15).Value = 64.3918918918919
$ws.Cells.Item(12, 16).Value = 37
$ws.Cells.Item(12, 17).Value = 6
$ws.Cells.Item(12, 18).Value = 6
$ws.Cells.Item(12, 19).Value = 6
$ws.Cells.Item(13, 6).Value = 347
$ws.Cells.Item(13, 12).Value = 'stimuli/img_rru0v.png'
$ws.Cells.Item(13, 13).Value = 56.45238095238095
$ws.Cells.Item(13, 14).Value = 39.42857142857143
$ws.Cells.Item(13, 15).Value = 47.94047619047619
$ws.Cells.Item(13, 16).Value = 42
$ws.Cells.Item(13, 17).Value = 4
$ws.Cells.Item(13, 18).Value = 4
$ws.Cells.Item(13, 19).Value = 4
$ws.Cells.Item(14, 6).Value = 348
$ws.Cells.Item(14, 12).Value = 'stimuli/img_9oofc.png'
$ws.Cells.Item(14, 13).Value = 82.47619047619048
$ws.Cells.Item(14, 14).Value = 65.5
$ws.Cells.Item(14, 15).Value = 73.98809523809524
$ws.Cells.Item(14, 16).Value = 42
$ws.Cells.Item(14, 17).Value = 8
$ws.Cells.Item(14, 18).Value = 8
$ws.Cells.Item(14, 19).Value = 8
$ws.Cells.Item(15, 6).Value = 349
$ws.Cells.Item(15, 8).Value = 'living_rooms'
$ws.Cells.Item(15, 9).Value = 'target'
$ws.Cells.Item(15, 11).Value = 'j'
$ws.Cells.Item(15, 12).Value = 'stimuli/img_tn8ys.png'
$ws.Cells.Item(15, 13).Value = 86.70454545454545
$ws.Cells.Item(15, 14).Value = 72.40909090909091
$ws.Cells.Item(15, 15).Value = 79.55681818181819
$ws.Cells.Item(15, 16).Value = 44
$ws.Cells.Item(15, 17).Value = 10
$ws.Cells.Item(15, 18).Value = 10
$ws.Cells.Item(15, 19).Value = 10
$ws.Cells.Item(16, 6).Value = 350
$ws.Cells.Item(16, 12).Value = 'stimuli/img_vh7v8.png'
$ws.Cells.Item(16, 13).Value = 78.70454545454545
$ws.Cells.Item(16, 14).Value = 59.63636363636363
$ws.Cells.Item(16, 15).Value = 69.17045454545455
$ws.Cells.Item(16, 16).Value = 44
$ws.Cells.Item(16, 17).Value = 7
$ws.Cells.Item(16, 18).Value = 7
$ws.Cells.Item(16, 19).Value = 7
$ws.Cells.Item(17, 6).Value = 351
$ws.Cells.Item(17, 12).Value = 'stimuli/img_ra2nm.png'
$ws.Cells.Item(17, 13).Value = 70.75
$ws.Cells.Item(17, 14).Value = 50.375
$ws.Cells.Item(17, 15).Value = 60.5625
$ws.Cells.Item(17, 16).Value = 40
$ws.Cells.Item(17, 17).Value = 6
$ws.Cells.Item(17, 18).Value = 6
$ws.Cells.Item(17, 19).Value = 6
$ws.Cells.Item(18, 6).Value = 352
$ws.Cells.Item(18, 8).Value = 'living_rooms'
$ws.Cells.Item(18, 9).Value = 'target'
$ws.Cells.Item(18, 11).Value = 'j'
$ws.Cells.Item(18, 12).Value = 'stimuli/img_vnxft.png'
$ws.Cells.Item(18, 13).Value = 53.22727272727273
$ws.Cells.Item(18, 14).Value = 34.84090909090909
$ws.Cells.Item(18, 15).Value = 44.03409090909091
$ws.Cells.Item(18, 16).Value = 44
$ws.Cells.Item(18, 17).Value = 3
$ws.Cells.Item(18, 18).Value = 3
$ws.Cells.Item(18, 19).Value = 3
$ws.Cells.Item(19, 6).Value = 353
$ws.Cells.Item(19, 12).Value = 'stimuli/img_5jp4f.png'
$ws.Cells.Item(19, 13).Value = 84.85714285714286
$ws.Cells.Item(19, 14).Value = 67.83333333333333
$ws.Cells.Item(19, 15).Value = 76.3452380952381
$ws.Cells.Item(19, 16).Value = 42
$ws.Cells.Item(19, 17).Value = 9
$ws.Cells.Item(19, 18).Value = 9
$ws.Cells.Item(19, 19).Value = 9
$ws.Cells.Item(20, 6).Value = 354
$ws.Cells.Item(20, 8).Value = 'living_rooms'
$ws.Cells.Item(20, 9).Value = 'target'
$ws.Cells.Item(20, 11).Value = 'j'
$ws.Cells.Item(20, 12).Value = 'stimuli/img_rg4in.png'
$ws.Cells.Item(20, 13).Value = 49.3695652173913
$ws.Cells.Item(20, 14).Value = 30.21739130434782
$ws.Cells.Item(20, 15).Value = 39.79347826086956
$ws.Cells.Item(20, 16).Value = 46
$ws.Cells.Item(20, 17).Value = 3
$ws.Cells.Item(20, 18).Value = 3
$ws.Cells.Item(20, 19).Value = 3
$ws.Cells.Item(21, 6).Value = 355
$ws.Cells.Item(21, 12).Value = 'stimuli/img_j4ttn.png'
$ws.Cells.Item(21, 13).Value = 12.61904761904762
$ws.Cells.Item(21, 14).Value = 11.42857142857143
$ws.Cells.Item(21, 15).Value = 12.02380952380952
$ws.Cells.Item(21, 16).Value = 42
$ws.Cells.Item(21, 17).Value = 1
$ws.Cells.Item(21, 18).Value = 1
$ws.Cells.Item(21, 19).Value = 1
$ws.Cells.Item(22, 6).Value = 356
$ws.Cells.Item(22, 12).Value = 'stimuli/img_swq34.png'
$ws.Cells.Item(22, 13).Value = 64.11363636363636
$ws.Cells.Item(22, 14).Value = 43.04545454545455
$ws.Cells.Item(22, 15).Value = 53.57954545454545
$ws.Cells.Item(22, 16).Value = 44
$ws.Cells.Item(22, 17).Value = 5
$ws.Cells.Item(22, 18).Value = 5
$ws.Cells.Item(22, 19).Value = 5
$ws.Cells.Item(23, 6).Value = 357
$ws.Cells.Item(23, 12).Value = 'stimuli/img_rych7.png'
$ws.Cells.Item(23, 13).Value = 30.4468085106383
$ws.Cells.Item(23, 14).Value = 23.4468085106383
$ws.Cells.Item(23, 15).Value = 26.9468085106383
$ws.Cells.Item(23, 16).Value = 47
$ws.Cells.Item(24, 6).Value = 358
$ws.Cells.Item(24, 12).Value = 'stimuli/img_njhlh.png'
$ws.Cells.Item(24, 13).Value = 59.74418604651163
$ws.Cells.Item(24, 14).Value = 41.51162790697674
$ws.Cells.Item(24, 15).Value = 50.62790697674419
$ws.Cells.Item(24, 16).Value = 43
$ws.Cells.Item(24, 17).Value = 4
$ws.Cells.Item(24, 18).Value = 4
$ws.Cells.Item(24, 19).Value = 4
$ws.Cells.Item(25, 6).Value = 359
$ws.Cells.Item(25, 8).Value = 'bedrooms'
$ws.Cells.Item(25, 9).Value = 'distractor'
$ws.Cells.Item(25, 11).Value = 'f'
$ws.Cells.Item(25, 12).Value = 'stimuli/img_jp28n.png'
$ws.Cells.Item(25, 13).Value = 65.02564102564102
$ws.Cells.Item(25, 14).Value = 44.97435897435897
$ws.Cells.Item(25, 15).Value = 55
$ws.Cells.Item(25, 16).Value = 39
$ws.Cells.Item(25, 17).Value = 4
$ws.Cells.Item(25, 18).Value = 4
$ws.Cells.Item(25, 19).Value = 4
$ws.Cells.Item(26, 6).Value = 360
$ws.Cells.Item(26, 12).Value = 'stimuli/img_of8d6.png'
$ws.Cells.Item(26, 13).Value = 26.04878048780488
$ws.Cells.Item(26, 14).Value = 19.14634146341463
$ws.Cells.Item(26, 15).Value = 22.59756097560975
$ws.Cells.Item(26, 16).Value = 41
$ws.Cells.Item(26, 17).Value = 1
$ws.Cells.Item(26, 18).Value = 1
$ws.Cells.Item(26, 19).Value = 1
$ws.Cells.Item(27, 6).Value = 361
$ws.Cells.Item(27, 8).Value = 'bedrooms'
$ws.Cells.Item(27, 9).Value = 'distractor'
$ws.Cells.Item(27, 11).Value = 'f'
$ws.Cells.Item(27, 12).Value = 'stimuli/img_d3t0o.png'
$ws.Cells.Item(27, 13).Value = 66.95121951219512
$ws.Cells.Item(27, 14).Value = 42.92682926829269
$ws.Cells.Item(27, 15).Value = 54.9390243902439
$ws.Cells.Item(27, 16).Value = 41
$ws.Cells.Item(27, 17).Value = 4
$ws.Cells.Item(27, 18).Value = 4
$ws.Cells.Item(27, 19).Value = 4
